$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-33 down to 24-34.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with data.
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44873
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 300000001
$ws.Cells.Item(23, 7).Value = "Rabanito"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 7900
$ws.Cells.Item(23, 11).Value = 3000
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 13).Value = 3000
$ws.Cells.Item(23, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(23, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(23, 16).Value = 30
$ws.Cells.Item(23, 17).Value = 100
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Match the date style (numFmtId 165) used by the rest of column D.
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
